$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 6) with a new Random-method sample.
$ws.Range("A6").Value = 42611.887638888889
$ws.Range("B6").Value = 31
$ws.Range("C6:M6").Value = 0
$ws.Range("N6").Value = "Random"
